# Bonus Drop Loot System Added
# Adds two new loot entries to the "Loot_All" sheet:
#   Row 13: BONUS_ARMOR_1 | Assets/Scripts/Weapons/ShotGun.prefab | Shot Gun | COMMON
#   Row 14: SHOT_GUN_1    | Assets/Scripts/Weapons/ShotGun.prefab | Shot Gun | COMMON

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loot_All")
$ws.Activate()

$ws.Range("A13").Value = "BONUS_ARMOR_1"
$ws.Range("B13").Value = "Assets/Scripts/Weapons/ShotGun.prefab"
$ws.Range("C13").Value = "Shot Gun"
$ws.Range("D13").Value = "COMMON"

$ws.Range("A14").Value = "SHOT_GUN_1"
$ws.Range("B14").Value = "Assets/Scripts/Weapons/ShotGun.prefab"
$ws.Range("C14").Value = "Shot Gun"
$ws.Range("D14").Value = "COMMON"

# Match the author's final selection state on the sheet
$ws.Range("C20").Select()
